# Re-sort the sighting rows: the underlying report rows were re-ordered,
# but each row's fixed site/date/observer metadata (columns C, I, J, O, P, S,
# T:Y, AA, AD:AG, AT, AW:AY) stays put on its original row number. Only the
# per-species identity columns move between rows:
#   A  Id
#   B  Taxonsorteringsordning
#   D  Rödlistade
#   E  TaxonId
#   F  Artnamn
#   G  Vetenskapligt namn
#   H  Auktor
#   K  Ålder-Stadium
#   L  Kön
#   M  Aktivitet
#   N  Metod
#   Q  Ost
#   R  Nord
#   AC Publik kommentar

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bands = @("A:B", "D:H", "K:N", "Q:R", "AC:AC")

function Get-RowBands($row) {
    $vals = @{}
    foreach ($band in $bands) {
        $parts = $band.Split(":")
        $addr = "$($parts[0])$row`:$($parts[1])$row"
        $vals[$band] = $ws.Range($addr).Value2
    }
    return $vals
}

function Set-RowBands($row, $vals) {
    foreach ($band in $bands) {
        $parts = $band.Split(":")
        $addr = "$($parts[0])$row`:$($parts[1])$row"
        $ws.Range($addr).Value2 = $vals[$band]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowBands $rowA
    $b = Get-RowBands $rowB
    Set-RowBands $rowA $b
    Set-RowBands $rowB $a
}

function Rotate-Rows([int[]]$rows) {
    # New-row(i) gets the content that used to live in Old-row(i+1),
    # wrapping the last back to the first: new[0]=old[1], new[1]=old[2], ..., new[last]=old[0]
    $originals = @()
    foreach ($r in $rows) { $originals += ,(Get-RowBands $r) }
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $src = $originals[($i + 1) % $rows.Length]
        Set-RowBands $rows[$i] $src
    }
}

# Pairwise full swaps
Swap-Rows 6 7
Swap-Rows 30 31
Swap-Rows 32 33

# Three-way rotations
Rotate-Rows @(15, 16, 17)
Rotate-Rows @(21, 23, 22)
